# Edit script for 棉布.xlsx
# 1) Remove columns F (棉布产销率) and G (棉布销售量) entirely, shifting remaining
#    cells left (dimension becomes A1:E65).
# 2) Within every year block of 4 quarterly rows (A,B,C,D), swap the B and C
#    quarter rows' full contents (label + data A:E) with each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete columns F:G, shifting cells left
$ws.Range("F1:G65").Delete()

# Step 2: swap B/C quarter rows for each year.
# Data starts at row 2 (2000年A) and is laid out in consecutive blocks of 4
# rows per year: row+0=A, row+1=B, row+2=C, row+3=D.
# Last data row is 65, so the last year block starts at row 62.
for ($yearStart = 2; $yearStart -le 62; $yearStart += 4) {
    $rowB = $yearStart + 1
    $rowC = $yearStart + 2

    $valuesB = $ws.Range("A" + $rowB + ":E" + $rowB).Value()
    $valuesC = $ws.Range("A" + $rowC + ":E" + $rowC).Value()

    $ws.Range("A" + $rowB + ":E" + $rowB).Value = $valuesC
    $ws.Range("A" + $rowC + ":E" + $rowC).Value = $valuesB
}
